$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C2" = "10-19"
    "D2" = "7-16"
    "G2" = "off"
    "B3" = "10-19"
    "D3" = "7-16"
    "G3" = "off"
    "B4" = "7-16"
    "D4" = "10-19"
    "F4" = "off"
    "G4" = "off"
    "B5" = "10-19"
    "D5" = "7-16"
    "G5" = "15-24"
    "C6" = "10-19"
    "D6" = "7-16"
    "G6" = "off"
    "D9" = "10-19"
    "F9" = "15-24"
    "C10" = "10-19"
    "F10" = "off"
    "F11" = "10-19"
    "G11" = "off"
    "C12" = "10-19"
    "G12" = "off"
    "B13" = "7-16"
    "B15" = "10-19"
    "G15" = "off"
    "B17" = "7-16"
    "D17" = "off"
    "B18" = "10-19"
    "G18" = "off"
    "B19" = "7-16"
    "D19" = "off"
    "F19" = "15-24"
    "E20" = "10-19"
    "F20" = "15-24"
    "C21" = "15-24"
    "F21" = "off"
    "B22" = "10-19"
    "C22" = "15-24"
    "F22" = "off"
    "B23" = "10-19"
    "D23" = "off"
    "F23" = "15-24"
    "B24" = "10-19"
    "F24" = "15-24"
    "B26" = "10-19"
    "D26" = "off"
    "B27" = "10-19"
    "D27" = "15-24"
    "G27" = "off"
    "B30" = "10-19"
    "D30" = "off"
    "E30" = "15-24"
    "B31" = "10-19"
    "D31" = "off"
    "E31" = "15-24"
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
